$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.712.25"
$ws.Range("E2").Value = "  -4.21%  "
$ws.Range("D3").Value = "1.816.82"
$ws.Range("E3").Value = "  -3.03%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'277.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.94%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.5079"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.28%  "
$ws.Range("E8").Value = "  -5.86%  "
$ws.Range("D9").Value = "'44.56"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("D10").Value = "'0.06659"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.31%  "
$ws.Range("D11").Value = "'20.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.40%  "
$ws.Range("D12").Value = "'0.8255"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.18%  "
$ws.Range("D13").Value = "'0.07874"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.41%  "
$ws.Range("D14").Value = "1.809.06"
$ws.Range("E14").Value = "  -3.39%  "
$ws.Range("D15").Value = "'5.072"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.46%  "
$ws.Range("D16").Value = "'87.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.45%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "'14.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.98%  "
$ws.Range("D19").Value = "'0.000008025"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.15%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "25.767.17"
$ws.Range("E21").Value = "  -4.13%  "
$ws.Range("D22").Value = "'4.739"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.02%  "
$ws.Range("D23").Value = "'9.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.12%  "
$ws.Range("D24").Value = "'6.095"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.78%  "
$ws.Range("D25").Value = "'142.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.69%  "
$ws.Range("D26").Value = "'2.207"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.66%  "
$ws.Range("D27").Value = "'1.676"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("E28").Value = "  -5.19%  "
$ws.Range("D29").Value = "'109.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.78%  "
$ws.Range("D30").Value = "'4.338"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.36%  "
$ws.Range("D31").Value = "'4.232"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.35%  "
$ws.Range("D32").Value = "'0.08794"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.86%  "
$ws.Range("D33").Value = "'0.04879"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.50%  "
$ws.Range("D34").Value = "'0.7287"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -10.77%  "
$ws.Range("D35").Value = "'1.136"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.22%  "
$ws.Range("D36").Value = "'2.873"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.76%  "
$ws.Range("D37").Value = "'1.000"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "'3.137"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.45%  "
$ws.Range("D39").Value = "'2.366"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.68%  "
$ws.Range("E40").Value = "  -5.30%  "
$ws.Range("D41").Value = "'0.5156"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -15.02%  "
$ws.Range("D42").Value = "'0.9640"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.95%  "
$ws.Range("D43").Value = "'6.220"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.27%  "
$ws.Range("D44").Value = "'110.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.98%  "
$ws.Range("D45").Value = "'8.025"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.05%  "
$ws.Range("D46").Value = "'1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'0.4560"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.85%  "
$ws.Range("D48").Value = "'0.1363"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.77%  "
$ws.Range("D49").Value = "'36.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.10%  "
$ws.Range("D50").Value = "'9.234"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.00%  "
$ws.Range("D51").Value = "'1.499"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.20%  "
